# V 2.0.2 se arreglo la fecha y hora de reimpresion
# Update patient record with the new (reprinted) data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Patient name and expediente number
$ws.Range("A6").Value = "DE LEON  SANTIZO  NEITAN  SAMUEL"
$ws.Range("G6").Value = "/201762107"

# Fecha de nacimiento y lugar de nacimiento
$ws.Range("A9").Value = "2016-09-02"
$ws.Range("E9").Value = "GUATEMALA"

# Ocupacion / Nacionalidad / Documento de identificacion -> now blank
$ws.Range("C11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("G11").Value = ""

# Emergency contact info
$ws.Range("A13").Value = "CLAUDIA SANTIZO"
$ws.Range("E13").Value = "16C. 19-10 ALAMEDA 2B Z.18"
$ws.Range("G13").Value = ""

# Fecha / hora de la asistencia medica (reimpresion)
$ws.Range("D14").Value = "Hora: 15:13:1"
$ws.Range("A15").Value = "24/10/2017"
